# Auto-generated edit script: refresh market-price derived columns (H-N)
# across multiple worksheets, per the scheduled-runner diff.
$wb = $excel.ActiveWorkbook

# ---------------- Sheet: ALC ----------------
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1548.75
$ws.Range("I4").Value = 117.8
$ws.Range("K4").Value = 117.8
$ws.Range("M4").Value = -3.799999999999997
$ws.Range("H17").Value = 2132.8518
$ws.Range("J17").Value = 2274.4583
$ws.Range("L17").Value = 6823.374899999999
$ws.Range("N17").Value = -7159.374899999999
$ws.Range("H40").Value = 4029.2593
$ws.Range("I40").Value = 3075.7368
$ws.Range("J40").Value = 6293.875
$ws.Range("K40").Value = 3075.7368
$ws.Range("L40").Value = 6293.875
$ws.Range("M40").Value = -2900.7368
$ws.Range("N40").Value = -6643.875
$ws.Range("H69").Value = 6909.033
$ws.Range("H70").Value = 4258.727
$ws.Range("I70").Value = 1400.5
$ws.Range("J70").Value = 4893.8887
$ws.Range("K70").Value = 4201.5
$ws.Range("L70").Value = 14681.6661
$ws.Range("M70").Value = -3931.5
$ws.Range("N70").Value = -15221.6661
$ws.Range("H72").Value = 6909.033
$ws.Range("H73").Value = 4258.727
$ws.Range("I73").Value = 1400.5
$ws.Range("J73").Value = 4893.8887
$ws.Range("K73").Value = 4201.5
$ws.Range("L73").Value = 14681.6661
$ws.Range("M73").Value = -3265.5
$ws.Range("N73").Value = -16553.6661
$ws.Range("H88").Value = 1525.8572
$ws.Range("I88").Value = 1619
$ws.Range("J88").Value = 1456
$ws.Range("K88").Value = 1619
$ws.Range("L88").Value = 1456
$ws.Range("M88").Value = -1213
$ws.Range("N88").Value = -2268
$ws.Range("H91").Value = 1525.8572
$ws.Range("I91").Value = 1619
$ws.Range("J91").Value = 1456
$ws.Range("K91").Value = 1619
$ws.Range("L91").Value = 1456
$ws.Range("M91").Value = -215
$ws.Range("N91").Value = -4264
$ws.Range("H112").Value = 1486.25
$ws.Range("I112").Value = 1486.25
$ws.Range("K112").Value = 4458.75
$ws.Range("M112").Value = -3350.75
$ws.Range("H137").Value = 2501.1177
$ws.Range("J137").Value = 3104.2222
$ws.Range("L137").Value = 9312.6666
$ws.Range("N137").Value = -14412.6666

# ---------------- Sheet: ARM ----------------
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("H32").Value = 3340929.8
$ws.Range("J32").Value = 14304214
$ws.Range("L32").Value = 14304214
$ws.Range("N32").Value = -14304788
$ws.Range("H43").Value = 10050000
$ws.Range("I43").Value = 20000000
$ws.Range("K43").Value = 20000000
$ws.Range("M43").Value = -19999687
$ws.Range("H74").Value = 4482.591
$ws.Range("I74").Value = 3798.0527
$ws.Range("K74").Value = 3798.0527
$ws.Range("M74").Value = -2924.0527
$ws.Range("H77").Value = 4482.591
$ws.Range("I77").Value = 3798.0527
$ws.Range("K77").Value = 18990.2635
$ws.Range("M77").Value = -14622.2635
$ws.Range("H88").Value = 772.4
$ws.Range("I88").Value = 462.5
$ws.Range("K88").Value = 462.5
$ws.Range("M88").Value = -56.5
$ws.Range("H91").Value = 772.4
$ws.Range("I91").Value = 462.5
$ws.Range("K91").Value = 462.5
$ws.Range("M91").Value = 941.5
$ws.Range("H135").Value = 60000
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 60000
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 60000
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -70140

# ---------------- Sheet: BSM ----------------
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 1293
$ws.Range("I36").Value = 1293
$ws.Range("K36").Value = 1293
$ws.Range("M36").Value = -759

# ---------------- Sheet: CRP ----------------
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9091.9375
$ws.Range("I31").Value = 3847.5
$ws.Range("J31").Value = 9841.143
$ws.Range("K31").Value = 3847.5
$ws.Range("L31").Value = 9841.143
$ws.Range("M31").Value = -3552.5
$ws.Range("N31").Value = -10431.143
$ws.Range("H34").Value = 9091.9375
$ws.Range("I34").Value = 3847.5
$ws.Range("J34").Value = 9841.143
$ws.Range("K34").Value = 3847.5
$ws.Range("L34").Value = 9841.143
$ws.Range("M34").Value = -3645.5
$ws.Range("N34").Value = -10245.143
$ws.Range("H58").Value = 6319.375
$ws.Range("I58").Value = 1305.5
$ws.Range("J58").Value = 7990.6665
$ws.Range("K58").Value = 1305.5
$ws.Range("L58").Value = 7990.6665
$ws.Range("M58").Value = -1102.5
$ws.Range("N58").Value = -8396.666499999999
$ws.Range("H105").Value = 1024.1111
$ws.Range("I105").Value = 900.75
$ws.Range("J105").Value = 2011
$ws.Range("K105").Value = 900.75
$ws.Range("L105").Value = 2011
$ws.Range("M105").Value = 846.25
$ws.Range("N105").Value = -5505
$ws.Range("H132").Value = 4827.9688
$ws.Range("I132").Value = 3997.6667
$ws.Range("K132").Value = 11993.0001
$ws.Range("M132").Value = -9463.000100000001
$ws.Range("H136").Value = 6319.375
$ws.Range("I136").Value = 1305.5
$ws.Range("J136").Value = 7990.6665
$ws.Range("K136").Value = 3916.5
$ws.Range("L136").Value = 23971.9995
$ws.Range("M136").Value = -1366.5
$ws.Range("N136").Value = -29071.9995

# ---------------- Sheet: CUL ----------------
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 4805.5
$ws.Range("I80").Value = 4684.1577
$ws.Range("J80").Value = 5266.6
$ws.Range("K80").Value = 14052.4731
$ws.Range("L80").Value = 15799.8
$ws.Range("M80").Value = -13116.4731
$ws.Range("N80").Value = -17671.8
$ws.Range("H81").Value = 1963.5
$ws.Range("I81").Value = 1719.5
$ws.Range("J81").Value = 2207.5
$ws.Range("K81").Value = 5158.5
$ws.Range("L81").Value = 6622.5
$ws.Range("M81").Value = -4035.5
$ws.Range("N81").Value = -8868.5
$ws.Range("H83").Value = 4805.5
$ws.Range("I83").Value = 4684.1577
$ws.Range("J83").Value = 5266.6
$ws.Range("K83").Value = 42157.41929999999
$ws.Range("L83").Value = 47399.4
$ws.Range("M83").Value = -37477.41929999999
$ws.Range("N83").Value = -56759.4
$ws.Range("H84").Value = 1963.5
$ws.Range("I84").Value = 1719.5
$ws.Range("J84").Value = 2207.5
$ws.Range("K84").Value = 15475.5
$ws.Range("L84").Value = 19867.5
$ws.Range("M84").Value = -9859.5
$ws.Range("N84").Value = -31099.5
$ws.Range("H121").Value = 1665.3334
$ws.Range("I121").Value = 190.8
$ws.Range("K121").Value = 572.4000000000001
$ws.Range("M121").Value = 737.5999999999999
$ws.Range("H138").Value = 5077.5713
$ws.Range("I138").Value = 2781.6667
$ws.Range("K138").Value = 8345.000100000001
$ws.Range("M138").Value = -3205.000100000001

# ---------------- Sheet: GSM ----------------
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5765.7144
$ws.Range("I70").Value = 2721.4
$ws.Range("J70").Value = 13376.5
$ws.Range("K70").Value = 2721.4
$ws.Range("L70").Value = 13376.5
$ws.Range("M70").Value = -2451.4
$ws.Range("N70").Value = -13916.5
$ws.Range("H73").Value = 5765.7144
$ws.Range("I73").Value = 2721.4
$ws.Range("J73").Value = 13376.5
$ws.Range("K73").Value = 2721.4
$ws.Range("L73").Value = 13376.5
$ws.Range("M73").Value = -1785.4
$ws.Range("N73").Value = -15248.5
$ws.Range("H132").Value = 47779.668
$ws.Range("I132").Value = 67923.56
$ws.Range("J132").Value = 7491.875
$ws.Range("K132").Value = 203770.68
$ws.Range("L132").Value = 22475.625
$ws.Range("M132").Value = -201240.68
$ws.Range("N132").Value = -27535.625

# ---------------- Sheet: LTW ----------------
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 1500
$ws.Range("J12").Value = 1500
$ws.Range("L12").Value = 1500
$ws.Range("N12").Value = -1840
$ws.Range("H46").Value = 5620.92
$ws.Range("I46").Value = 4692.8
$ws.Range("J46").Value = 6239.6665
$ws.Range("K46").Value = 4692.8
$ws.Range("L46").Value = 6239.6665
$ws.Range("M46").Value = -4504.8
$ws.Range("N46").Value = -6615.6665
$ws.Range("H63").Value = 44444
$ws.Range("I63").Value = 44444
$ws.Range("K63").Value = 44444
$ws.Range("M63").Value = -43695
$ws.Range("H66").Value = 44444
$ws.Range("I66").Value = 44444
$ws.Range("K66").Value = 133332
$ws.Range("M66").Value = -129588
$ws.Range("H68").Value = 9284.9
$ws.Range("I68").Value = 8524.75
$ws.Range("K68").Value = 8524.75
$ws.Range("M68").Value = -7775.75
$ws.Range("H71").Value = 9284.9
$ws.Range("I71").Value = 8524.75
$ws.Range("K71").Value = 42623.75
$ws.Range("M71").Value = -38879.75
$ws.Range("H82").Value = 3316.25
$ws.Range("I82").Value = 1419.2
$ws.Range("K82").Value = 1419.2
$ws.Range("M82").Value = -1058.2
$ws.Range("H85").Value = 3316.25
$ws.Range("I85").Value = 1419.2
$ws.Range("K85").Value = 1419.2
$ws.Range("M85").Value = -171.2
$ws.Range("H136").Value = 2807
$ws.Range("I136").Value = 2020.125
$ws.Range("K136").Value = 6060.375
$ws.Range("M136").Value = -3510.375

# ---------------- Sheet: WVR ----------------
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 943.3333
$ws.Range("I122").Value = 943.3333
$ws.Range("K122").Value = 2829.9999
$ws.Range("M122").Value = -379.9998999999998
$ws.Range("H136").Value = 2101
$ws.Range("I136").Value = 1311.6
$ws.Range("J136").Value = 9995
$ws.Range("K136").Value = 3934.8
$ws.Range("L136").Value = 29985
$ws.Range("M136").Value = -1384.8
$ws.Range("N136").Value = -35085
